$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append two new data rows (21 and 22) that duplicate rows 1 and 2,
# but belong to "section" 12 instead of section 1.
$ws.Range("A21").Value = $ws.Range("A1").Value2
$ws.Range("B21").Value = $ws.Range("B1").Value2
$ws.Range("C21").Value = $ws.Range("C1").Value2
$ws.Range("D21").Value = 12
$ws.Range("E21").Value = 1

$ws.Range("A22").Value = $ws.Range("A2").Value2
$ws.Range("B22").Value = $ws.Range("B2").Value2
$ws.Range("C22").Value = $ws.Range("C2").Value2
$ws.Range("D22").Value = 12
$ws.Range("E22").Value = 2

# Keep the printed page in portrait orientation (as saved in the file).
$ws.PageSetup.Orientation = 1

# Match the selection left behind in the saved file.
$ws.Range("G25").Select()
